$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'27.549.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.39%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'1.622.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.69%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'  +0.32%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'212.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.10%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'0.521"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.85%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("E7").Value = "'  +0.32%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'22.92"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.46%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'0.260"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.75%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("E10").Value = "'  +0.26%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.0890"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.41%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'1.852.62"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.71%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'1.631.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.32%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'4.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.53%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'0.548"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.43%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'64.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.36%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'27.592.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.70%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'229.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.68%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'0.0₃0721"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.28%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'7.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.57%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("E21").Value = "'  +0.20%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'4.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.26%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'9.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.52%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'2.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +5.34%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'149.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.41%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'6.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.27%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("B27").Value = "'BinanceUSD"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.15%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("B28").Value = "'Stellar"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'0.111"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.10%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'15.56"
$ws.Range("D29").Style = "Normal"

# Row 30
$ws.Range("E30").Value = "'  -0.32%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("E31").Value = "'  -0.87%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("E32").Value = "'  +0.15%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'1.457.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +3.55%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'3.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -2.98%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'1.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.49%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("E36").Value = "'  +0.37%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("E37").Value = "'  -0.51%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("E38").Value = "'  +0.46%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'0.865"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.70%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.907"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.12%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("E41").Value = "'  +7.49%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("E42").Value = "'  +0.26%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("E43").Value = "'  -1.69%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'2.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.30%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'5.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.74%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'2.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.42%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'1.763.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.70%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("E48").Value = "'  +1.51%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'86.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.19%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'0.0₆0100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -5.60%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("B51").Value = "'EnergySwap"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'7.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.48%  "
$ws.Range("E51").Style = "Normal"
